# "Export with no is_pref and no lev distance"
#
# The data rows (2-10) get re-ordered/relabelled: the `id` (col B) and
# `speaker_variant` (col C) pairs are rewritten to a new arrangement, and
# the `is_prefered` column (col D) - which used to hold "x" markers - is
# cleared out entirely for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (id, speaker_variant) pairs for rows 2..10, in order.
$data = @(
    @("#lammert",  "Lammert"),
    @("#griet",    "Griet"),
    @("#grietjen", "Grietjen"),
    @("#schout",   "Schout"),
    @("#trvn",     "Trvn"),
    @("#moer",     "Moer"),
    @("#dienaer",  "Dienaer"),
    @("#buurwyf",  "Buurwyf"),
    @("#tryn",     "Tryn")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]   # column B: id
    $ws.Cells.Item($row, 3).Value = $data[$i][1]   # column C: speaker_variant
    $ws.Cells.Item($row, 4).Value = ""              # column D: is_prefered -> cleared
}
